$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B7:AW7").Value = "Configurations/SSHConfig.csv"
